$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Description")

# Row 8 ("Notes"): trim the trailing "Please see SCRGSP (2016) for more information." sentence
$ws.Range("B8").Value = "The scope of this indicator has been increased to include all public hospitals reporting to the Non-Admitted Patient Emergency Department Care National Minimum Data Set. Data for 2013" + [char]8211 + "14 have been resupplied for the revised scope, but it is not possible to provide comparable data for the years prior to 2013" + [char]8211 + "14."
$ws.Rows.Item(8).RowHeight = 30.55

# Row 9: add "Source" label in column A, replace column B text with the new citation
$ws.Range("A9").Value = "Source"
$ws.Range("B9").Value = "AIHW (various years), Australian hospital statistics, Health Services Series"
$ws.Rows.Item(9).RowHeight = 13.8

# New row 10: second source citation line, matching B9's wrapped text formatting
$ws.Range("B10").Value = "AIHW (various years), Emergency department care: Australian hospital statistics, Health services series."
$ws.Range("B10").WrapText = $true
$ws.Range("B10").Font.Name = $ws.Range("B9").Font.Name
$ws.Range("B10").Font.Size = $ws.Range("B9").Font.Size
$ws.Rows.Item(10).RowHeight = 25.45
